$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix two standalone quantity corrections (rows 34 and 36) ---
$ws.Range("C34").Value = 0
$ws.Range("C36").Value = 0

# --- Rows 40-89: Remessa (A, forced text), Material (B), Quantidade (C) ---
$matA = New-Object "object[,]" 50,1
$matB = New-Object "object[,]" 50,1
$matC = New-Object "object[,]" 50,1

$matA[0,0] = '=TEXT(80266697,"0")'
$matB[0,0] = "15403-TDK-N"
$matC[0,0] = 180
$matA[1,0] = '=TEXT(80266697,"0")'
$matB[1,0] = "15403-TDK-N"
$matC[1,0] = 20
$matA[2,0] = '=TEXT(80266697,"0")'
$matB[2,0] = "40343-TDK-N"
$matC[2,0] = 1130
$matA[3,0] = '=TEXT(80266697,"0")'
$matB[3,0] = "40343-TDK-N"
$matC[3,0] = 670
$matA[4,0] = '=TEXT(80266700,"0")'
$matB[4,0] = "10361-ARI-I"
$matC[4,0] = 1
$matA[5,0] = '=TEXT(80266701,"0")'
$matB[5,0] = "10255-ARI-I"
$matC[5,0] = 1
$matA[6,0] = '=TEXT(80266702,"0")'
$matB[6,0] = "10399-ARI-I"
$matC[6,0] = 1
$matA[7,0] = '=TEXT(80266703,"0")'
$matB[7,0] = "10650-ARI-I"
$matC[7,0] = 1
$matA[8,0] = '=TEXT(80266704,"0")'
$matB[8,0] = "20041-CTY-I"
$matC[8,0] = 1
$matA[9,0] = '=TEXT(80266706,"0")'
$matB[9,0] = "21340-NZX-I"
$matC[9,0] = 1
$matA[10,0] = '=TEXT(80266707,"0")'
$matB[10,0] = "10377-ARI-I"
$matC[10,0] = 1
$matA[11,0] = '=TEXT(80266708,"0")'
$matB[11,0] = "10252-ARI-I"
$matC[11,0] = 1
$matA[12,0] = '=TEXT(80266709,"0")'
$matB[12,0] = "20935-CTY-I"
$matC[12,0] = 1
$matA[13,0] = '=TEXT(80266710,"0")'
$matB[13,0] = "10000-LDG-I"
$matC[13,0] = 1
$matA[14,0] = '=TEXT(80266711,"0")'
$matB[14,0] = "30159-OSR-I"
$matC[14,0] = 72000
$matA[15,0] = '=TEXT(80266712,"0")'
$matB[15,0] = "10382-ARI-I"
$matC[15,0] = 1
$matA[16,0] = '=TEXT(80266713,"0")'
$matB[16,0] = "20988-CTY-I"
$matC[16,0] = 1
$matA[17,0] = '=TEXT(80266714,"0")'
$matB[17,0] = "10020-ARI-I"
$matC[17,0] = 1
$matA[18,0] = '=TEXT(80266715,"0")'
$matB[18,0] = "10253-ARI-I"
$matC[18,0] = 1
$matA[19,0] = '=TEXT(80266716,"0")'
$matB[19,0] = "10001-LDG-I"
$matC[19,0] = 1
$matA[20,0] = '=TEXT(80266719,"0")'
$matB[20,0] = "15188-DLO-I"
$matC[20,0] = 30
$matA[21,0] = '=TEXT(80266719,"0")'
$matB[21,0] = "15187-DLO-I"
$matC[21,0] = 16
$matA[22,0] = '=TEXT(80266719,"0")'
$matB[22,0] = "15124-DLO-I"
$matC[22,0] = 12
$matA[23,0] = '=TEXT(80266719,"0")'
$matB[23,0] = "15330-DLO-I"
$matC[23,0] = 5
$matA[24,0] = '=TEXT(80266719,"0")'
$matB[24,0] = "15260-DLO-I"
$matC[24,0] = 10
$matA[25,0] = '=TEXT(80266719,"0")'
$matB[25,0] = "10533-DLO-I"
$matC[25,0] = 1
$matA[26,0] = '=TEXT(80266719,"0")'
$matB[26,0] = "11619-DLO-I"
$matC[26,0] = 2
$matA[27,0] = '=TEXT(80266719,"0")'
$matB[27,0] = "14763-DLO-I"
$matC[27,0] = 5
$matA[28,0] = '=TEXT(80266719,"0")'
$matB[28,0] = "15141-DLO-I"
$matC[28,0] = 15
$matA[29,0] = '=TEXT(80266719,"0")'
$matB[29,0] = "15178-DLO-I"
$matC[29,0] = 1
$matA[30,0] = '=TEXT(80266719,"0")'
$matB[30,0] = "15259-DLO-I"
$matC[30,0] = 5
$matA[31,0] = '=TEXT(80266720,"0")'
$matB[31,0] = "10369-ARI-I"
$matC[31,0] = 1
$matA[32,0] = '=TEXT(80266721,"0")'
$matB[32,0] = "10645-ARI-I"
$matC[32,0] = 1
$matA[33,0] = '=TEXT(80266722,"0")'
$matB[33,0] = "60339-YAG-I"
$matC[33,0] = 20
$matA[34,0] = '=TEXT(80266723,"0")'
$matB[34,0] = "14122-TDK-N"
$matC[34,0] = 500
$matA[35,0] = '=TEXT(80266723,"0")'
$matB[35,0] = "40322-TDK-N"
$matC[35,0] = 500
$matA[36,0] = '=TEXT(80266723,"0")'
$matB[36,0] = "30029-KMT-I"
$matC[36,0] = 500
$matA[37,0] = '=TEXT(80266723,"0")'
$matB[37,0] = "60162-OUT-L"
$matC[37,0] = 500
$matA[38,0] = '=TEXT(80266723,"0")'
$matB[38,0] = "22793-STM-I"
$matC[38,0] = 300
$matA[39,0] = '=TEXT(80266723,"0")'
$matB[39,0] = "17741-SET-I"
$matC[39,0] = 500
$matA[40,0] = '=TEXT(80266724,"0")'
$matB[40,0] = "10288-ARI-I"
$matC[40,0] = 1
$matA[41,0] = '=TEXT(80266725,"0")'
$matB[41,0] = "10453-ARI-I"
$matC[41,0] = 1
$matA[42,0] = '=TEXT(80266726,"0")'
$matB[42,0] = "10334-ARI-I"
$matC[42,0] = 1
$matA[43,0] = '=TEXT(80266727,"0")'
$matB[43,0] = "30100-ROY-I"
$matC[43,0] = 4000
$matA[44,0] = '=TEXT(80266727,"0")'
$matB[44,0] = "10802-ROY-I"
$matC[44,0] = 35000
$matA[45,0] = '=TEXT(80266727,"0")'
$matB[45,0] = "20637-TDK-I"
$matC[45,0] = 500
$matA[46,0] = '=TEXT(80266729,"0")'
$matB[46,0] = "30152-OSR-I"
$matC[46,0] = 50000
$matA[47,0] = '=TEXT(80266729,"0")'
$matB[47,0] = "30173-OSR-L"
$matC[47,0] = 120000
$matA[48,0] = '=TEXT(80266731,"0")'
$matB[48,0] = "10689-ARI-I"
$matC[48,0] = 1
$matA[49,0] = '=TEXT(80266732,"0")'
$matB[49,0] = "10130-ARI-I"
$matC[49,0] = 1

$rngA = $ws.Range("A40:A89")
$rngB = $ws.Range("B40:B89")
$rngC = $ws.Range("C40:C89")

# Column A: write as formulas first so pure-digit Remessa codes stay text
# (TEXT() forces string type), then freeze to static values via paste-special.
$rngA.Formula = $matA
$rngA.Copy() | Out-Null
$rngA.PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

$rngB.Value = $matB
$rngC.Value = $matC

# --- Selection / view state (A2:C89 per the saved workbook) ---
$ws.Range("A2:C89").Select() | Out-Null
